$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3251533742331288
$ws.Range("C2").Value = 0.2269938650306748
$ws.Range("B3").Value = 0.2392638036809816
$ws.Range("C3").Value = 0.2085889570552147
